$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Cells.Item(2, 4).Value = 61004
$ws.Cells.Item(2, 5).Value = 649
$ws.Cells.Item(2, 6).Value = 17
$ws.Cells.Item(2, 7).Value = 8162
$ws.Cells.Item(2, 8).Value = 5089
$ws.Cells.Item(2, 9).Value = 5027
$ws.Cells.Item(2, 10).Value = 62
$ws.Cells.Item(2, 11).Value = 77186
$ws.Cells.Item(2, 12).Value = 30758
$ws.Cells.Item(2, 13).Value = 46428
$ws.Cells.Item(2, 14).Value = 45534
$ws.Cells.Item(2, 15).Value = 894
$ws.Cells.Item(2, 16).Value = 3880
$ws.Cells.Item(2, 17).Value = 3726
$ws.Cells.Item(2, 18).Value = -3722
$ws.Cells.Item(2, 19).Value = -218
$ws.Cells.Item(2, 20).Value = 8599
$ws.Cells.Item(2, 21).Value = -4872
$ws.Cells.Item(2, 22).Value = 17131
$ws.Cells.Item(2, 23).Value = 1.06
$ws.Cells.Item(2, 24).Value = 8.34
$ws.Cells.Item(2, 25).Value = 11.52
$ws.Cells.Item(2, 26).Value = 6.83
$ws.Cells.Item(2, 27).Value = 66.25
$ws.Cells.Item(2, 28).Value = 914.89
$ws.Cells.Item(2, 29).Value = 6478
$ws.Cells.Item(2, 30).Value = 8.44
$ws.Cells.Item(2, 31).Value = 58718
$ws.Cells.Item(2, 32).Value = 0.93
$ws.Cells.Item(2, 33).Value = 750
$ws.Cells.Item(2, 34).Value = 1.37
$ws.Cells.Item(2, 35).Value = 11.6
$ws.Cells.Item(2, 36).Value = 74693696

# Row 3
$ws.Cells.Item(3, 4).Value = 61763
$ws.Cells.Item(3, 5).Value = 3013
$ws.Cells.Item(3, 6).Value = 3013
$ws.Cells.Item(3, 7).Value = 3668
$ws.Cells.Item(3, 8).Value = 206
$ws.Cells.Item(3, 9).Value = 112
$ws.Cells.Item(3, 10).Value = 95
$ws.Cells.Item(3, 11).Value = 72695
$ws.Cells.Item(3, 12).Value = 29541
$ws.Cells.Item(3, 13).Value = 43154
$ws.Cells.Item(3, 14).Value = 42221
$ws.Cells.Item(3, 15).Value = 933
$ws.Cells.Item(3, 16).Value = 3880
$ws.Cells.Item(3, 17).Value = 5371
$ws.Cells.Item(3, 18).Value = -2212
$ws.Cells.Item(3, 19).Value = 303
$ws.Cells.Item(3, 20).Value = 11962
$ws.Cells.Item(3, 21).Value = -6591
$ws.Cells.Item(3, 22).Value = 20416
$ws.Cells.Item(3, 23).Value = 4.88
$ws.Cells.Item(3, 24).Value = 0.33
$ws.Cells.Item(3, 25).Value = 0.26
$ws.Cells.Item(3, 26).Value = 0.27
$ws.Cells.Item(3, 27).Value = 68.45
$ws.Cells.Item(3, 28).Value = 899.67
$ws.Cells.Item(3, 29).Value = 144
$ws.Cells.Item(3, 30).Value = 436.26
$ws.Cells.Item(3, 31).Value = 55887
$ws.Cells.Item(3, 32).Value = 1.13
$ws.Cells.Item(3, 33).Value = 500
$ws.Cells.Item(3, 34).Value = 0.79
$ws.Cells.Item(3, 35).Value = 338.89
$ws.Cells.Item(3, 36).Value = 74693696

# Row 4
$ws.Cells.Item(4, 4).Value = 60330
$ws.Cells.Item(4, 5).Value = 244
$ws.Cells.Item(4, 6).Value = 244
$ws.Cells.Item(4, 7).Value = 321
$ws.Cells.Item(4, 8).Value = 229
$ws.Cells.Item(4, 9).Value = 147
$ws.Cells.Item(4, 10).Value = 82
$ws.Cells.Item(4, 11).Value = 76626
$ws.Cells.Item(4, 12).Value = 33250
$ws.Cells.Item(4, 13).Value = 43376
$ws.Cells.Item(4, 14).Value = 42401
$ws.Cells.Item(4, 15).Value = 975
$ws.Cells.Item(4, 16).Value = 3880
$ws.Cells.Item(4, 17).Value = 6796
$ws.Cells.Item(4, 18).Value = -11863
$ws.Cells.Item(4, 19).Value = 2813
$ws.Cells.Item(4, 20).Value = 10519
$ws.Cells.Item(4, 21).Value = -3722
$ws.Cells.Item(4, 22).Value = 24436
$ws.Cells.Item(4, 23).Value = 0.4
$ws.Cells.Item(4, 24).Value = 0.38
$ws.Cells.Item(4, 25).Value = 0.35
$ws.Cells.Item(4, 26).Value = 0.31
$ws.Cells.Item(4, 27).Value = 76.66
$ws.Cells.Item(4, 28).Value = 922.4
$ws.Cells.Item(4, 29).Value = 190
$ws.Cells.Item(4, 30).Value = 268.04
$ws.Cells.Item(4, 31).Value = 56125
$ws.Cells.Item(4, 32).Value = 0.91
$ws.Cells.Item(4, 33).Value = 500
$ws.Cells.Item(4, 34).Value = 0.98
$ws.Cells.Item(4, 35).Value = 257.8
$ws.Cells.Item(4, 36).Value = 74693696

# Row 5
$ws.Cells.Item(5, 4).Value = 68385
$ws.Cells.Item(5, 5).Value = 3062
$ws.Cells.Item(5, 6).Value = 3062
$ws.Cells.Item(5, 7).Value = 2535
$ws.Cells.Item(5, 8).Value = 1773
$ws.Cells.Item(5, 9).Value = 1617
$ws.Cells.Item(5, 10).Value = 155
$ws.Cells.Item(5, 11).Value = 77674
$ws.Cells.Item(5, 12).Value = 34359
$ws.Cells.Item(5, 13).Value = 43315
$ws.Cells.Item(5, 14).Value = 42316
$ws.Cells.Item(5, 15).Value = 998
$ws.Cells.Item(5, 16).Value = 3880
$ws.Cells.Item(5, 17).Value = 7177
$ws.Cells.Item(5, 18).Value = -12323
$ws.Cells.Item(5, 19).Value = 1961
$ws.Cells.Item(5, 20).Value = 14763
$ws.Cells.Item(5, 21).Value = -7586
$ws.Cells.Item(5, 22).Value = 25689
$ws.Cells.Item(5, 23).Value = 4.48
$ws.Cells.Item(5, 24).Value = 2.59
$ws.Cells.Item(5, 25).Value = 3.82
$ws.Cells.Item(5, 26).Value = 2.3
$ws.Cells.Item(5, 27).Value = 79.32
$ws.Cells.Item(5, 28).Value = 942.01
$ws.Cells.Item(5, 29).Value = 2084
$ws.Cells.Item(5, 30).Value = 47.98
$ws.Cells.Item(5, 31).Value = 56013
$ws.Cells.Item(5, 32).Value = 1.79
$ws.Cells.Item(5, 33).Value = 750
$ws.Cells.Item(5, 34).Value = 0.75
$ws.Cells.Item(5, 35).Value = 35.12
$ws.Cells.Item(5, 36).Value = 74693696

# Row 6
$ws.Cells.Item(6, 4).Value = 81930
$ws.Cells.Item(6, 5).Value = 10181
$ws.Cells.Item(6, 6).Value = 10181
$ws.Cells.Item(6, 7).Value = 9404
$ws.Cells.Item(6, 8).Value = 6850
$ws.Cells.Item(6, 9).Value = 6562
$ws.Cells.Item(6, 11).Value = 86449
$ws.Cells.Item(6, 12).Value = 36984
$ws.Cells.Item(6, 13).Value = 49465
$ws.Cells.Item(6, 14).Value = 48210
$ws.Cells.Item(6, 16).Value = 3880
$ws.Cells.Item(6, 17).Value = 15587
$ws.Cells.Item(6, 18).Value = -6981
$ws.Cells.Item(6, 19).Value = -2358
$ws.Cells.Item(6, 20).Value = 11952
$ws.Cells.Item(6, 21).Value = 3635
$ws.Cells.Item(6, 22).Value = 24541
$ws.Cells.Item(6, 23).Value = 12.43
$ws.Cells.Item(6, 24).Value = 8.36
$ws.Cells.Item(6, 25).Value = 14.5
$ws.Cells.Item(6, 26).Value = 8.35
$ws.Cells.Item(6, 27).Value = 74.77
$ws.Cells.Item(6, 28).Value = 1088.79
$ws.Cells.Item(6, 29).Value = 8457
$ws.Cells.Item(6, 30).Value = 12.24
$ws.Cells.Item(6, 31).Value = 63815
$ws.Cells.Item(6, 32).Value = 1.62
$ws.Cells.Item(6, 33).Value = 1000
$ws.Cells.Item(6, 34).Value = 0.97
$ws.Cells.Item(6, 35).Value = 11.53
$ws.Cells.Item(6, 36).Value = 74693696

# Row 7
$ws.Cells.Item(7, 4).Value = 83348
$ws.Cells.Item(7, 5).Value = 6864
$ws.Cells.Item(7, 7).Value = 6420
$ws.Cells.Item(7, 8).Value = 5708
$ws.Cells.Item(7, 9).Value = 5609
$ws.Cells.Item(7, 11).Value = 88751
$ws.Cells.Item(7, 12).Value = 33848
$ws.Cells.Item(7, 13).Value = 54903
$ws.Cells.Item(7, 14).Value = 53448
$ws.Cells.Item(7, 16).Value = 3880
$ws.Cells.Item(7, 17).Value = 13599
$ws.Cells.Item(7, 18).Value = -5412
$ws.Cells.Item(7, 19).Value = -4812
$ws.Cells.Item(7, 20).Value = 11518
$ws.Cells.Item(7, 21).Value = 3298
$ws.Cells.Item(7, 23).Value = 8.24
$ws.Cells.Item(7, 24).Value = 6.85
$ws.Cells.Item(7, 25).Value = 11.04
$ws.Cells.Item(7, 26).Value = 6.52
$ws.Cells.Item(7, 27).Value = 61.65
$ws.Cells.Item(7, 29).Value = 7228
$ws.Cells.Item(7, 30).Value = 18.61
$ws.Cells.Item(7, 31).Value = 70748
$ws.Cells.Item(7, 32).Value = 1.9
$ws.Cells.Item(7, 33).Value = 972
$ws.Cells.Item(7, 34).Value = 0.72
$ws.Cells.Item(7, 35).Value = 12.95

# Row 8
$ws.Cells.Item(8, 4).Value = 85717
$ws.Cells.Item(8, 5).Value = 8303
$ws.Cells.Item(8, 7).Value = 8048
$ws.Cells.Item(8, 8).Value = 6010
$ws.Cells.Item(8, 9).Value = 5791
$ws.Cells.Item(8, 11).Value = 94562
$ws.Cells.Item(8, 12).Value = 34496
$ws.Cells.Item(8, 13).Value = 60066
$ws.Cells.Item(8, 14).Value = 58614
$ws.Cells.Item(8, 16).Value = 3880
$ws.Cells.Item(8, 17).Value = 13840
$ws.Cells.Item(8, 18).Value = -10118
$ws.Cells.Item(8, 19).Value = -1652
$ws.Cells.Item(8, 20).Value = 9952
$ws.Cells.Item(8, 21).Value = 4248
$ws.Cells.Item(8, 23).Value = 9.69
$ws.Cells.Item(8, 24).Value = 7.01
$ws.Cells.Item(8, 25).Value = 10.37
$ws.Cells.Item(8, 26).Value = 6.54
$ws.Cells.Item(8, 27).Value = 57.43
$ws.Cells.Item(8, 29).Value = 7463
$ws.Cells.Item(8, 30).Value = 16.75
$ws.Cells.Item(8, 31).Value = 77586
$ws.Cells.Item(8, 32).Value = 1.61
$ws.Cells.Item(8, 33).Value = 1022
$ws.Cells.Item(8, 34).Value = 0.82
$ws.Cells.Item(8, 35).Value = 13.18

# Row 9
$ws.Cells.Item(9, 4).Value = 92708
$ws.Cells.Item(9, 5).Value = 9839
$ws.Cells.Item(9, 7).Value = 9482
$ws.Cells.Item(9, 8).Value = 7117
$ws.Cells.Item(9, 9).Value = 6878
$ws.Cells.Item(9, 11).Value = 100972
$ws.Cells.Item(9, 12).Value = 34234
$ws.Cells.Item(9, 13).Value = 66738
$ws.Cells.Item(9, 14).Value = 64926
$ws.Cells.Item(9, 16).Value = 3880
$ws.Cells.Item(9, 17).Value = 15511
$ws.Cells.Item(9, 18).Value = -10735
$ws.Cells.Item(9, 19).Value = -1484
$ws.Cells.Item(9, 20).Value = 9921
$ws.Cells.Item(9, 21).Value = 5194
$ws.Cells.Item(9, 23).Value = 10.61
$ws.Cells.Item(9, 24).Value = 7.68
$ws.Cells.Item(9, 25).Value = 11.13
$ws.Cells.Item(9, 26).Value = 7.28
$ws.Cells.Item(9, 27).Value = 51.3
$ws.Cells.Item(9, 29).Value = 8863
$ws.Cells.Item(9, 30).Value = 14.1
$ws.Cells.Item(9, 31).Value = 85940
$ws.Cells.Item(9, 32).Value = 1.45
$ws.Cells.Item(9, 33).Value = 1064
$ws.Cells.Item(9, 34).Value = 0.85
$ws.Cells.Item(9, 35).Value = 11.55
